$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column E (Date Sampled), shifting E:N to H:Q
$ws.Range("E1:G1").EntireColumn.Insert()

# Set header labels for the new columns
$ws.Range("E1").Value = "Day"
$ws.Range("F1").Value = "Month"
$ws.Range("G1").Value = "Year"

# Fill Day/Month/Year values derived from the Date Sampled column (now column H)
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = 2015

$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 2015

$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = 2015

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 2015

# Update selection to match target
$ws.Range("E4").Select()

Write-Output "done"
